$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.572.09"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.513.99"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.11"
$ws.Range("E5").Value = "  +4.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.58"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  +2.41%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.47"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.75"
$ws.Range("E12").Value = "  +3.62%  "

$ws.Range("E13").Value = "  -2.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.898.18"
$ws.Range("E14").Value = "  -1.46%  "

$ws.Range("E15").Value = "  +6.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.522.88"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.856"
$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.548.17"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  -2.93%  "

$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.51"
$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.57"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.06"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.01"
$ws.Range("E26").Value = "  -2.42%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  +11.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.11"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.59"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("E31").Value = "  -0.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.34"
$ws.Range("E32").Value = "  +0.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.24"
$ws.Range("E33").Value = "  +5.26%  "

$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("E37").Value = "  -4.29%  "

$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("E39").Value = "  +1.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.76"
$ws.Range("E40").Value = "  -7.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.39"
$ws.Range("E41").Value = "  +1.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.84"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.02"
$ws.Range("E43").Value = "  -2.91%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.018.30"
$ws.Range("E46").Value = "  -2.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.36"
$ws.Range("E47").Value = "  -4.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.95"
$ws.Range("E48").Value = "  -2.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.755.43"
$ws.Range("E49").Value = "  -1.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.02"
$ws.Range("E50").Value = "  -1.84%  "

$ws.Range("E51").Value = "  +0.96%  "
